$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241101-085714-"

# Shift the reference date (column G) forward by one day (2024-10-31 -> 2024-11-01) for every data row (2-274)
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value2 = 45597
}

# Update projected balance (D), forecast balance (E) and total balance (H) for the rows that changed
$ws.Cells.Item(5, 5).Value2 = 1013.56
$ws.Cells.Item(5, 8).Value2 = 1013.56

$ws.Cells.Item(6, 5).Value2 = 105.23
$ws.Cells.Item(6, 8).Value2 = 105.23

$ws.Cells.Item(8, 5).Value2 = 50.73
$ws.Cells.Item(8, 8).Value2 = 50.73

$ws.Cells.Item(15, 5).Value2 = 972.75
$ws.Cells.Item(15, 8).Value2 = 972.75

$ws.Cells.Item(17, 5).Value2 = 152.26
$ws.Cells.Item(17, 8).Value2 = 152.26

$ws.Cells.Item(43, 5).Value2 = 465.95
$ws.Cells.Item(43, 8).Value2 = 465.95

$ws.Cells.Item(49, 5).Value2 = 304.4
$ws.Cells.Item(49, 8).Value2 = 304.4

$ws.Cells.Item(51, 5).Value2 = 11765.61
$ws.Cells.Item(51, 8).Value2 = 11765.61

$ws.Cells.Item(52, 5).Value2 = 344.91
$ws.Cells.Item(52, 8).Value2 = 344.91

$ws.Cells.Item(57, 5).Value2 = 423.53
$ws.Cells.Item(57, 8).Value2 = 423.53

$ws.Cells.Item(60, 5).Value2 = 1013.11
$ws.Cells.Item(60, 8).Value2 = 1013.11

$ws.Cells.Item(97, 5).Value2 = 76.14
$ws.Cells.Item(97, 8).Value2 = 76.14

$ws.Cells.Item(99, 5).Value2 = 1444.6
$ws.Cells.Item(99, 8).Value2 = 1444.6

$ws.Cells.Item(101, 5).Value2 = 152.24
$ws.Cells.Item(101, 8).Value2 = 152.24

$ws.Cells.Item(102, 5).Value2 = 228.68
$ws.Cells.Item(102, 8).Value2 = 228.68

$ws.Cells.Item(103, 5).Value2 = 1001.47
$ws.Cells.Item(103, 8).Value2 = 1001.47

$ws.Cells.Item(104, 5).Value2 = 1034.51
$ws.Cells.Item(104, 8).Value2 = 1034.51

$ws.Cells.Item(105, 5).Value2 = 845.87
$ws.Cells.Item(105, 8).Value2 = 845.87

$ws.Cells.Item(107, 5).Value2 = 22323.97
$ws.Cells.Item(107, 8).Value2 = 22323.97

$ws.Cells.Item(108, 4).Value2 = 0
$ws.Cells.Item(108, 5).Value2 = -136.49
$ws.Cells.Item(108, 8).Value2 = -136.49

$ws.Cells.Item(109, 5).Value2 = 0.19
$ws.Cells.Item(109, 8).Value2 = 0.19

$ws.Cells.Item(110, 5).Value2 = 897.36
$ws.Cells.Item(110, 8).Value2 = 897.36

$ws.Cells.Item(112, 5).Value2 = 2.73
$ws.Cells.Item(112, 8).Value2 = 2.73

$ws.Cells.Item(113, 5).Value2 = 3.28
$ws.Cells.Item(113, 8).Value2 = 3.28

$ws.Cells.Item(114, 5).Value2 = 0
$ws.Cells.Item(114, 8).Value2 = 0

$ws.Cells.Item(118, 5).Value2 = 828.9
$ws.Cells.Item(118, 8).Value2 = 828.9

$ws.Cells.Item(132, 5).Value2 = 1005.3
$ws.Cells.Item(132, 8).Value2 = 1005.3

$ws.Cells.Item(138, 5).Value2 = 1862.15
$ws.Cells.Item(138, 8).Value2 = 1862.15

$ws.Cells.Item(143, 4).Value2 = -18596.22
$ws.Cells.Item(143, 5).Value2 = 73747.61
$ws.Cells.Item(143, 8).Value2 = 55151.39

$ws.Cells.Item(148, 5).Value2 = 0.04
$ws.Cells.Item(148, 8).Value2 = 0.04

$ws.Cells.Item(158, 5).Value2 = 297.47
$ws.Cells.Item(158, 8).Value2 = 297.47

$ws.Cells.Item(165, 4).Value2 = 0
$ws.Cells.Item(165, 5).Value2 = 33.13
$ws.Cells.Item(165, 8).Value2 = 33.13

$ws.Cells.Item(173, 5).Value2 = 1020.51
$ws.Cells.Item(173, 8).Value2 = 1020.51

$ws.Cells.Item(189, 4).Value2 = 0
$ws.Cells.Item(189, 5).Value2 = 0.12
$ws.Cells.Item(189, 8).Value2 = 0.12

$ws.Cells.Item(224, 5).Value2 = 560.11
$ws.Cells.Item(224, 8).Value2 = 560.11

$ws.Cells.Item(230, 5).Value2 = 6397.75
$ws.Cells.Item(230, 8).Value2 = 6397.75

$ws.Cells.Item(232, 5).Value2 = 33881.82
$ws.Cells.Item(232, 8).Value2 = 33881.82

$ws.Cells.Item(235, 5).Value2 = 964.44
$ws.Cells.Item(235, 8).Value2 = 964.44

$ws.Cells.Item(249, 5).Value2 = 776.7
$ws.Cells.Item(249, 8).Value2 = 776.7

$ws.Cells.Item(255, 5).Value2 = 684.91
$ws.Cells.Item(255, 8).Value2 = 684.91

$ws.Cells.Item(264, 5).Value2 = 964.68
$ws.Cells.Item(264, 8).Value2 = 964.68

$ws.Cells.Item(265, 5).Value2 = 1013.79
$ws.Cells.Item(265, 8).Value2 = 1013.79

$ws.Cells.Item(270, 5).Value2 = 1009.92
$ws.Cells.Item(270, 8).Value2 = 1009.92

$ws.Cells.Item(271, 5).Value2 = 1013.63
$ws.Cells.Item(271, 8).Value2 = 1013.63

$ws.Cells.Item(273, 5).Value2 = 1008.9
$ws.Cells.Item(273, 8).Value2 = 1008.9
